$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number (e.g. "553.83") need to be
# forced to Text format first, otherwise Excel would store them as numeric
# values instead of literal strings. We reset the style back to Normal
# afterwards so no visible formatting change remains on the cell.
$textCells = @("D5", "D6", "D10", "D12", "D14", "D17", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D28", "D29", "D30", "D34", "D38", "D39", "D40", "D48", "D49", "D50")
foreach ($ref in $textCells) {
    $ws.Range($ref).NumberFormat = "@"
}

$ws.Range("D2").Value = "64.342.84"
$ws.Range("E2").Value = "  +0.40%  "
$ws.Range("D3").Value = "3.330.51"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "553.83"
$ws.Range("E5").Value = "  +0.60%  "
$ws.Range("D6").Value = "173.61"
$ws.Range("E6").Value = "  +0.63%  "
$ws.Range("E7").Value = "  +2.02%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "3.318.77"
$ws.Range("E9").Value = "  -0.18%  "
$ws.Range("D10").Value = "0.172"
$ws.Range("E10").Value = "  +6.88%  "
$ws.Range("E11").Value = "  +2.11%  "
$ws.Range("D12").Value = "53.33"
$ws.Range("E12").Value = "  +0.33%  "
$ws.Range("E13").Value = "  +2.38%  "
$ws.Range("D14").Value = "9.09"
$ws.Range("E14").Value = "  +0.97%  "
$ws.Range("D15").Value = "3.862.98"
$ws.Range("E15").Value = "  +0.15%  "
$ws.Range("E16").Value = "  +3.05%  "
$ws.Range("D17").Value = "18.16"
$ws.Range("E17").Value = "  -0.53%  "
$ws.Range("D18").Value = "3.341.31"
$ws.Range("E18").Value = "  +0.06%  "
$ws.Range("B19").Value = "WrappedBTC"
$ws.Range("C19").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D19").Value = "64.381.97"
$ws.Range("E19").Value = "  +0.67%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").Value = "11.76"
$ws.Range("E20").Value = "  +0.24%  "
$ws.Range("D21").Value = "0.988"
$ws.Range("E21").Value = "  +1.74%  "
$ws.Range("D22").Value = "452.83"
$ws.Range("E22").Value = "  +6.76%  "
$ws.Range("D23").Value = "4.98"
$ws.Range("E23").Value = "  +5.85%  "
$ws.Range("D24").Value = "4.06"
$ws.Range("E24").Value = "  -0.23%  "
$ws.Range("D25").Value = "88.03"
$ws.Range("E25").Value = "  +4.82%  "
$ws.Range("D26").Value = "13.88"
$ws.Range("E26").Value = "  +5.09%  "
$ws.Range("E27").Value = "  +2.63%  "
$ws.Range("D28").Value = "10.58"
$ws.Range("E28").Value = "  -0.33%  "
$ws.Range("B29").Value = "Filecoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D29").Value = "8.60"
$ws.Range("E29").Value = "  +0.34%  "
$ws.Range("B30").Value = "EthereumClassic"
$ws.Range("C30").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D30").Value = "30.96"
$ws.Range("E30").Value = "  +4.47%  "
$ws.Range("E31").Value = "  -2.27%  "
$ws.Range("E32").Value = "  +0.37%  "
$ws.Range("E33").Value = "  +6.57%  "
$ws.Range("D34").Value = "570.19"
$ws.Range("E34").Value = "  -3.99%  "
$ws.Range("E35").Value = "  +0.37%  "
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("E37").Value = "  -0.22%  "
$ws.Range("D38").Value = "3.52"
$ws.Range("E38").Value = "  +0.89%  "
$ws.Range("D39").Value = "35.44"
$ws.Range("E39").Value = "  +0.27%  "
$ws.Range("D40").Value = "0.367"
$ws.Range("E40").Value = "  +0.89%  "
$ws.Range("D41").Value = "0.0₃0732"
$ws.Range("E41").Value = "  -2.19%  "
$ws.Range("D42").Value = "3.067.34"
$ws.Range("E42").Value = "  -1.05%  "
$ws.Range("E43").Value = "  +2.61%  "
$ws.Range("E44").Value = "  -1.14%  "
$ws.Range("E45").Value = "  -0.01%  "
$ws.Range("E46").Value = "  +0.65%  "
$ws.Range("E47").Value = "  +3.95%  "
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").Value = "142.54"
$ws.Range("E48").Value = "  +7.19%  "
$ws.Range("B49").Value = "FirstDigitalUSD"
$ws.Range("C49").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D49").Value = "1.00"
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("D50").Value = "2.51"
$ws.Range("E50").Value = "  -2.46%  "
$ws.Range("E51").Value = "  -0.17%  "

foreach ($ref in $textCells) {
    $ws.Range($ref).Style = "Normal"
}
